# Add a new worksheet "tryeditorcode" after "invalidcode", matching the
# look/feel of the existing "validcode" sheet (same header/box formatting),
# and make it the active (selected) sheet.

$wb = $excel.ActiveWorkbook

$validcode = $wb.Worksheets.Item("validcode")
$invalidcode = $wb.Worksheets.Item("invalidcode")

# Insert the new sheet right after "invalidcode" (i.e. as the new last sheet).
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $invalidcode)
$newSheet.Name = "tryeditorcode"

# Copy the cell formatting (fill/border styles) of the small "validcode"
# table (header row + one data row) onto our new sheet, then reuse the
# data-row formatting for the extra third row.
$validcode.Range("A1:B2").Copy()
$newSheet.Range("A1:B2").PasteSpecial(-4122)
$validcode.Range("A2").Copy()
$newSheet.Range("A3").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the values.
$newSheet.Range("A1").Value = "pythonCode"
$newSheet.Range("B1").Value = "output"
$newSheet.Range("A2").Value = 'print("Hello")'
$newSheet.Range("B2").Value = "Hello"
$newSheet.Range("A3").Value = "Hello"

# Match the recorded selection on the new (now active) sheet.
[void]$newSheet.Range("A1:B4").Select()
